# Fixed naive component forecaster bug - Presentation state 11.02.
# Update forecast-error table values in rows 2-11 (columns B:G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = -0.06304482099587455; "C2" = 1.284473961807544;  "D2" = 3.268916597491724;  "E2" = 1.808014545708005;  "F2" = 1.824894733127932;  "G2" = 51
    "B3" = 0.3848683813723368;   "C3" = 1.606026431035804;  "D3" = 7.828542112970349;  "E3" = 2.797953200639773;  "F3" = 2.799493110457711;  "G3" = 50
    "B4" = 0.8746662378248481;   "C4" = 1.947065274388025;  "D4" = 10.77140143740715;  "E4" = 3.281981328010132;  "F4" = 3.196064181409324;  "G4" = 49
    "B5" = 0.5898522373503049;   "C5" = 1.91023936529586;   "D5" = 11.52381436931002;  "E5" = 3.394674412857589;  "F5" = 3.378412877893673;  "G5" = 48
    "B6" = 0.6796283455196681;   "C6" = 2.053046112058658;  "D6" = 12.25639243636365;  "E6" = 3.500913086091062;  "F6" = 3.471440604892087;  "G6" = 47
    "B7" = 0.4620040488475629;   "C7" = 1.878470610035362;  "D7" = 11.8468855781181;   "E7" = 3.441930501639757;  "F7" = 3.456566949417224;  "G7" = 38
    "B8" = 0.5682316542479313;   "C8" = 1.945817767357757;  "D8" = 12.30934837299113;  "E8" = 3.508468094908536;  "F8" = 3.50990290671253;   "G8" = 37
    "B9" = 0.5226750354899881;   "C9" = 2.444751333587228;  "D9" = 20.24450866695628;  "E9" = 4.499389810513897;  "F9" = 4.585023615187054;  "G9" = 20
    "B10" = -0.5732758996450275; "C10" = 2.122935804733599; "D10" = 8.802873317106101; "E10" = 2.96696365281176;  "F10" = 3.02991975663685;  "G10" = 13
    "B11" = -0.3393998147840875; "C11" = 2.457134838492891; "D11" = 7.456238556435196; "E11" = 2.730611388761718; "F11" = 3.029242133389085; "G11" = 5
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
